# Update the yearly database: shift each fiscal-year column one period to the
# left (1396->1397, 1397->1398, 1398->1399, 1399->1400, 1400->1401) and bring
# in the newly reported last column (1401/12), per the "update database and
# change read_price algorithm" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column headers (rows 8 and 24): shift the reporting-period labels ----
$ws.Range("E8").Value2  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value2  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value2  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value2  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value2  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value2 = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value2 = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value2 = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value2 = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value2 = "دوازده ماهه منتهی به 1401/12"

# ---- هزینه حمل و نقل و انتقال (row 10) ----
$ws.Range("E10").Value2 = 101998
$ws.Range("F10").Value2 = 207693
$ws.Range("G10").Value2 = 458796
$ws.Range("H10").Value2 = 796569
$ws.Range("I10").Value2 = 1234661

# ---- حق العمل و کمیسیون فروش (row 13) ----
$ws.Range("E13").Value2 = 281927
$ws.Range("F13").Value2 = 291034
$ws.Range("G13").Value2 = 483308
$ws.Range("H13").Value2 = 784530
$ws.Range("I13").Value2 = 1234661

# ---- هزینه تبلیغات (row 14) ----
$ws.Range("E14").Value2 = 8770
$ws.Range("F14").Value2 = 12328
$ws.Range("G14").Value2 = 31194
$ws.Range("H14").Value2 = 51822
$ws.Range("I14").Value2 = 21598

# ---- هزینه مواد مصرفی (row 15) ----
$ws.Range("E15").Value2 = 1052
$ws.Range("F15").Value2 = 1684
$ws.Range("G15").Value2 = 2498
$ws.Range("H15").Value2 = 3198
$ws.Range("I15").Value2 = 5298

# ---- هزینه انرژی (آب، برق، گاز و سوخت) (row 16) ----
$ws.Range("E16").Value2 = 4306
$ws.Range("F16").Value2 = 4536
$ws.Range("G16").Value2 = 9221
$ws.Range("H16").Value2 = 23269
$ws.Range("I16").Value2 = 25344

# ---- هزینه استهلاک (row 17) ----
$ws.Range("E17").Value2 = 90946
$ws.Range("F17").Value2 = 131033
$ws.Range("G17").Value2 = 209007
$ws.Range("H17").Value2 = 326384
$ws.Range("I17").Value2 = 544924

# ---- هزینه حقوق و دستمزد (row 18) ----
$ws.Range("E18").Value2 = 3940
$ws.Range("F18").Value2 = 0
$ws.Range("G18").Value2 = 0
$ws.Range("H18").Value2 = 0
$ws.Range("I18").Value2 = 0

# ---- سایر هزینه ها (row 19) ----
$ws.Range("E19").Value2 = 110676
$ws.Range("F19").Value2 = 73591
$ws.Range("G19").Value2 = 154853
$ws.Range("H19").Value2 = 274974
$ws.Range("I19").Value2 = 268894

# ---- جمع (row 20) ----
$ws.Range("E20").Value2 = 603615
$ws.Range("F20").Value2 = 721899
$ws.Range("G20").Value2 = 1348877
$ws.Range("H20").Value2 = 2260746
$ws.Range("I20").Value2 = 3335380

# ---- تعداد پرسنل غیر تولیدی شرکت (row 26) ----
$ws.Range("E26").Value2 = 141
$ws.Range("F26").Value2 = 179
$ws.Range("G26").Value2 = 212
$ws.Range("H26").Value2 = 741
$ws.Range("I26").Value2 = 741

# ---- تعداد پرسنل تولیدی شرکت (row 27) ----
$ws.Range("E27").Value2 = 295
$ws.Range("F27").Value2 = 322
$ws.Range("G27").Value2 = 550
$ws.Range("H27").Value2 = 240
$ws.Range("I27").Value2 = 240

$wb.Save()
